# Update LR-pair TPM-derived statistics (Tnc-Itgav) with refreshed NATMI
# output values (new TPM run). Only the numeric value cells change; the
# text columns (A-D) and cell counts (K, L) are unaffected.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 6.449754000000001
$ws.Range("H2").Value = 19.349262
$ws.Range("I2").Value = 0.03479900749229446
$ws.Range("J2").Value = 0.03479900749229446
$ws.Range("M2").Value = 9.423852333333334
$ws.Range("N2").Value = 28.271557
$ws.Range("O2").Value = 0.06654336290212845
$ws.Range("P2").Value = 0.06654336290212845
$ws.Range("Q2").Value = 60.78152928232602
$ws.Range("R2").Value = 547.0337635409342
$ws.Range("S2").Value = 0.002315642984193637
$ws.Range("T2").Value = 0.002315642984193637

$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 6.449754000000001
$ws.Range("H3").Value = 19.349262
$ws.Range("I3").Value = 0.03479900749229446
$ws.Range("J3").Value = 0.03479900749229446
$ws.Range("O3").Value = 0.3572423751649123
$ws.Range("P3").Value = 0.3572423751649123
$ws.Range("Q3").Value = 326.3095963290921
$ws.Range("R3").Value = 2936.786366961828
$ws.Range("S3").Value = 0.01243168008992885
$ws.Range("T3").Value = 0.01243168008992885

$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 6.449754000000001
$ws.Range("H4").Value = 19.349262
$ws.Range("I4").Value = 0.03479900749229446
$ws.Range("J4").Value = 0.03479900749229446
$ws.Range("M4").Value = 26.84076266666667
$ws.Range("N4").Value = 80.522288
$ws.Range("O4").Value = 0.1895270158659356
$ws.Range("P4").Value = 0.1895270158659356
$ws.Range("Q4").Value = 173.116316372384
$ws.Range("R4").Value = 1558.046847351456
$ws.Range("S4").Value = 0.006595352045110903
$ws.Range("T4").Value = 0.006595352045110904

$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 6.449754000000001
$ws.Range("H5").Value = 19.349262
$ws.Range("I5").Value = 0.03479900749229446
$ws.Range("J5").Value = 0.03479900749229446
$ws.Range("M5").Value = 54.762539
$ws.Range("N5").Value = 164.287617
$ws.Range("O5").Value = 0.3866872460670236
$ws.Range("P5").Value = 0.3866872460670236
$ws.Range("Q5").Value = 353.2049049654061
$ws.Range("R5").Value = 3178.844144688654
$ws.Range("S5").Value = 0.01345633237306107
$ws.Range("T5").Value = 0.01345633237306107

$ws.Range("I6").Value = 0.663783921437469
$ws.Range("J6").Value = 0.6637839214374691
$ws.Range("M6").Value = 9.423852333333334
$ws.Range("N6").Value = 28.271557
$ws.Range("O6").Value = 0.06654336290212845
$ws.Range("P6").Value = 0.06654336290212845
$ws.Range("Q6").Value = 1159.39518869676
$ws.Range("R6").Value = 10434.55669827084
$ws.Range("S6").Value = 0.04417041437281142
$ws.Range("T6").Value = 0.04417041437281143

$ws.Range("I7").Value = 0.663783921437469
$ws.Range("J7").Value = 0.6637839214374691
$ws.Range("O7").Value = 0.3572423751649123
$ws.Range("P7").Value = 0.3572423751649123
$ws.Range("Q7").Value = 6224.2885375959
$ws.Range("R7").Value = 56018.5968383631
$ws.Range("S7").Value = 0.237131744690601
$ws.Range("T7").Value = 0.237131744690601

$ws.Range("I8").Value = 0.663783921437469
$ws.Range("J8").Value = 0.6637839214374691
$ws.Range("M8").Value = 26.84076266666667
$ws.Range("N8").Value = 80.522288
$ws.Range("O8").Value = 0.1895270158659356
$ws.Range("P8").Value = 0.1895270158659356
$ws.Range("Q8").Value = 3302.158182871035
$ws.Range("R8").Value = 29719.42364583931
$ws.Range("S8").Value = 0.1258049858098321
$ws.Range("T8").Value = 0.1258049858098322

$ws.Range("I9").Value = 0.663783921437469
$ws.Range("J9").Value = 0.6637839214374691
$ws.Range("M9").Value = 54.762539
$ws.Range("N9").Value = 164.287617
$ws.Range("O9").Value = 0.3866872460670236
$ws.Range("P9").Value = 0.3866872460670236
$ws.Range("Q9").Value = 6737.311026494086
$ws.Range("R9").Value = 60635.79923844677
$ws.Range("S9").Value = 0.2566767765642245
$ws.Range("T9").Value = 0.2566767765642245

$ws.Range("G10").Value = 55.79038633333334
$ws.Range("H10").Value = 167.371159
$ws.Range("I10").Value = 0.3010114916028843
$ws.Range("J10").Value = 0.3010114916028843
$ws.Range("M10").Value = 9.423852333333334
$ws.Range("N10").Value = 28.271557
$ws.Range("O10").Value = 0.06654336290212845
$ws.Range("P10").Value = 0.06654336290212845
$ws.Range("Q10").Value = 525.7603624249515
$ws.Range("R10").Value = 4731.843261824563
$ws.Range("S10").Value = 0.02003031692344172
$ws.Range("T10").Value = 0.02003031692344172

$ws.Range("G11").Value = 55.79038633333334
$ws.Range("H11").Value = 167.371159
$ws.Range("I11").Value = 0.3010114916028843
$ws.Range("J11").Value = 0.3010114916028843
$ws.Range("O11").Value = 0.3572423751649123
$ws.Range("P11").Value = 0.3572423751649123
$ws.Range("Q11").Value = 2822.578728347483
$ws.Range("R11").Value = 25403.20855512734
$ws.Range("S11").Value = 0.1075340602121474
$ws.Range("T11").Value = 0.1075340602121474

$ws.Range("G12").Value = 55.79038633333334
$ws.Range("H12").Value = 167.371159
$ws.Range("I12").Value = 0.3010114916028843
$ws.Range("J12").Value = 0.3010114916028843
$ws.Range("M12").Value = 26.84076266666667
$ws.Range("N12").Value = 80.522288
$ws.Range("O12").Value = 0.1895270158659356
$ws.Range("P12").Value = 0.1895270158659356
$ws.Range("Q12").Value = 1497.456518654644
$ws.Range("R12").Value = 13477.10866789179
$ws.Range("S12").Value = 0.05704980974484877
$ws.Range("T12").Value = 0.05704980974484879

$ws.Range("G13").Value = 55.79038633333334
$ws.Range("H13").Value = 167.371159
$ws.Range("I13").Value = 0.3010114916028843
$ws.Range("J13").Value = 0.3010114916028843
$ws.Range("M13").Value = 54.762539
$ws.Range("N13").Value = 164.287617
$ws.Range("O13").Value = 0.3866872460670236
$ws.Range("P13").Value = 0.3866872460670236
$ws.Range("Q13").Value = 3055.223207404234
$ws.Range("R13").Value = 27497.0088666381
$ws.Range("S13").Value = 0.1163973047224463
$ws.Range("T13").Value = 0.1163973047224463

$ws.Range("E14").Value = 2
$ws.Range("F14").Value = 0.6666666666666666
$ws.Range("G14").Value = 0.07517133333333333
$ws.Range("H14").Value = 0.225514
$ws.Range("I14").Value = 0.0004055794673521549
$ws.Range("J14").Value = 0.000405579467352155
$ws.Range("M14").Value = 9.423852333333334
$ws.Range("N14").Value = 28.271557
$ws.Range("O14").Value = 0.06654336290212845
$ws.Range("P14").Value = 0.06654336290212845
$ws.Range("Q14").Value = 0.7084035450331111
$ws.Range("R14").Value = 6.375631905298
$ws.Range("S14").Value = 0.0000269886216816664
$ws.Range("T14").Value = 0.00002698862168166641

$ws.Range("E15").Value = 2
$ws.Range("F15").Value = 0.6666666666666666
$ws.Range("G15").Value = 0.07517133333333333
$ws.Range("H15").Value = 0.225514
$ws.Range("I15").Value = 0.0004055794673521549
$ws.Range("J15").Value = 0.000405579467352155
$ws.Range("O15").Value = 0.3572423751649123
$ws.Range("P15").Value = 0.3572423751649123
$ws.Range("Q15").Value = 3.803110542746222
$ws.Range("R15").Value = 34.227994884716
$ws.Range("S15").Value = 0.0001448901722350038
$ws.Range("T15").Value = 0.0001448901722350039

$ws.Range("E16").Value = 2
$ws.Range("F16").Value = 0.6666666666666666
$ws.Range("G16").Value = 0.07517133333333333
$ws.Range("H16").Value = 0.225514
$ws.Range("I16").Value = 0.0004055794673521549
$ws.Range("J16").Value = 0.000405579467352155
$ws.Range("M16").Value = 26.84076266666667
$ws.Range("N16").Value = 80.522288
$ws.Range("O16").Value = 0.1895270158659356
$ws.Range("P16").Value = 0.1895270158659356
$ws.Range("Q16").Value = 2.017655917336889
$ws.Range("R16").Value = 18.158903256032
$ws.Range("S16").Value = 0.00007686826614374956
$ws.Range("T16").Value = 0.00007686826614374959

$ws.Range("E17").Value = 2
$ws.Range("F17").Value = 0.6666666666666666
$ws.Range("G17").Value = 0.07517133333333333
$ws.Range("H17").Value = 0.225514
$ws.Range("I17").Value = 0.0004055794673521549
$ws.Range("J17").Value = 0.000405579467352155
$ws.Range("M17").Value = 54.762539
$ws.Range("N17").Value = 164.287617
$ws.Range("O17").Value = 0.3866872460670236
$ws.Range("P17").Value = 0.3866872460670236
$ws.Range("Q17").Value = 4.116573073348666
$ws.Range("R17").Value = 37.04915766013799
$ws.Range("S17").Value = 0.0001568324072917351
$ws.Range("T17").Value = 0.0001568324072917351
